# Rows 48-53 of the sheet describe separate fungal records that were
# re-synced from source data: several rows were re-keyed onto a different
# underlying record, picking up that record's identity / name / coordinate
# columns, while the per-record "time" columns (Starttid / Sluttid, Z & AB)
# were dropped entirely. Columns A,B,D,E,F,G,H,P,Q,R cycle between rows
# 48->53, the rest of each row is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move as a block per row.
$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")

# Snapshot the "before" values for every row in the affected block so that
# writing the new values for one row never clobbers data still needed for
# another (this permutes, it does not just overwrite with new constants).
$snapshot = @{}
foreach ($r in 48..53) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# after-row -> before-row the block of values is sourced from.
$mapping = @{
    48 = 53
    49 = 52
    50 = 51
    51 = 49
    52 = 48
    53 = 50
}

foreach ($r in 48..53) {
    $src = $snapshot[$mapping[$r]]
    foreach ($col in $cols) {
        $val = $src[$col]
        if ($col -eq "Q" -or $col -eq "R") {
            $val = [Math]::Round([double]$val)
        }
        $ws.Range("$col$r").Value = $val
    }
    # Starttid (Z) and Sluttid (AB) are removed for these rows.
    $ws.Range("Z$r").ClearContents()
    $ws.Range("AB$r").ClearContents()
}
